# Update the "Förändrad" (changed) date column for rows 2-7 from 45208 to 45212.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C7").Value = 45212
